$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.279.35'
$ws.Range("D3").Value = '1.866.54'
$ws.Range("E3").Value = '  +0.23%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.68'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.70%  '
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2855'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.82%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06571'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.38'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07828'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.32%  '
$ws.Range("B12").Value = 'Litecoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '96.88'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.89%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.847.50'
$ws.Range("E13").Value = '  -0.80%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6983'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.53%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.094'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.03%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '268.05'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.56%  '
$ws.Range("D17").Value = '30.427.22'
$ws.Range("E17").Value = '  +0.56%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.84'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.55%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007642'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.00%  '
$ws.Range("D21").Value = '2.139.26'
$ws.Range("E21").Value = '  +1.66%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.234'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.71%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.175'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.455'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.78%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '166.79'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.36%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.87'
$ws.Range("D27").Style = "Normal"
$ws.Range("E28").Value = '  -0.76%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.368'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.92%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09910'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.60%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.355'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.04%  '
$ws.Range("E32").Value = '  -0.94%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.046'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.30%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04724'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.17%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.131'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.26%  '
$ws.Range("E36").Value = '  +0.33%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.717'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01873'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.16%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.753'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.05%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.334'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.07%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '72.92'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.40%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.949'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.14%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4171'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.34%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.000'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.07%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8363'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.07%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '103.15'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.09%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '971.96'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.96%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.114'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.62%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.104'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.25%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.47'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.89%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05682'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.36%  '